# Duplicate the "Map Chart" slide (slide 10) and retitle the copy.
$p = $ppt.ActivePresentation
$src = $p.Slides.Item(10)

$range = $src.Duplicate()
$copy = $range.Item(1)

foreach ($shp in $copy.Shapes) {
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Map Chart") {
            $shp.TextFrame.TextRange.Text = "Map Chart - copy"
        }
    }
}

# Materialize a notes page for the new slide (mirrors the source slide having
# a notesSlide part); the host only supports the notes body placeholder.
$notes = $copy.NotesPage
$notes.Shapes.AddPlaceholder(2) | Out-Null
